{"js": "// Update the date heading and the 5x5 grid of two-digit multiplication\n// answers (stored in the single table, one populated row every 5 rows).\nconst body = context.document.body;\n\n// --- date paragraph -------------------------------------------------\nconst dateHits = body.search(\"2025-01-03 Friday\", { matchCase: true });\ndateHits.load(\"items\");\nawait context.sync();\nif (dateHits.items.length > 0) {\n  dateHits.items[0].insertText(\"2025-01-04 Saturday\", \"Replace\");\n}\n\n// --- table of answers -------------------------------------------------\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Row indices (0-based) that hold data: 0, 4, 9, 14, 19 \u2014 each followed\n// by 4 blank spacer rows. New values listed left-to-right, top-to-bottom.\nconst newValues = [\n  [\"49\u00d736=1764\", \"36\u00d790=3240\", \"16\u00d733=528\", \"37\u00d712=444\", \"36\u00d781=2916\"],\n  [\"95\u00d770=6650\", \"11\u00d785=935\", \"57\u00d780=4560\", \"85\u00d727=2295\", \"93\u00d754=5022\"],\n  [\"25\u00d798=2450\", \"46\u00d734=1564\", \"45\u00d713=585\", \"57\u00d732=1824\", \"77\u00d727=2079\"],\n  [\"40\u00d760=2400\", \"22\u00d798=2156\", \"96\u00d740=3840\", \"99\u00d735=3465\", \"16\u00d785=1360\"],\n  [\"60\u00d749=2940\", \"39\u00d749=1911\", \"69\u00d721=1449\", \"70\u00d789=6230\", \"56\u00d780=4480\"],\n];\nconst dataRowIndexes = [0, 4, 9, 14, 19];\n\nconst cellsCollections = [];\nfor (const rowIdx of dataRowIndexes) {\n  const cells = rows.items[rowIdx].cells;\n  cells.load(\"items\");\n  cellsCollections.push(cells);\n}\nawait context.sync();\n\nfor (let i = 0; i < dataRowIndexes.length; i++) {\n  const cells = cellsCollections[i].items;\n  const vals = newValues[i];\n  for (let c = 0; c < vals.length; c++) {\n    cells[c].value = vals[c];\n  }\n}\nawait context.sync();\n", "ps1": "# Update the date heading and the 5x5 grid of two-digit multiplication\n# answers (stored in the single table, one populated row every 5 rows).\n$d = $word.ActiveDocument\n\n# --- date paragraph -------------------------------------------------\n$d.Paragraphs.Item(1).Range.Text = \"2025-01-04 Saturday\"\n\n# --- table of answers -------------------------------------------------\n$t = $d.Tables.Item(1)\n\n# Row numbers (1-based) that hold data: 1, 5, 10, 15, 20 \u2014 each followed\n# by 4 blank spacer rows. New values listed left-to-right, top-to-bottom.\n$dataRows = @(1, 5, 10, 15, 20)\n$newValues = @(\n    @(\"49\u00d736=1764\", \"36\u00d790=3240\", \"16\u00d733=528\", \"37\u00d712=444\", \"36\u00d781=2916\"),\n    @(\"95\u00d770=6650\", \"11\u00d785=935\", \"57\u00d780=4560\", \"85\u00d727=2295\", \"93\u00d754=5022\"),\n    @(\"25\u00d798=2450\", \"46\u00d734=1564\", \"45\u00d713=585\", \"57\u00d732=1824\", \"77\u00d727=2079\"),\n    @(\"40\u00d760=2400\", \"22\u00d798=2156\", \"96\u00d740=3840\", \"99\u00d735=3465\", \"16\u00d785=1360\"),\n    @(\"60\u00d749=2940\", \"39\u00d749=1911\", \"69\u00d721=1449\", \"70\u00d789=6230\", \"56\u00d780=4480\")\n)\n\nfor ($r = 0; $r -lt $dataRows.Length; $r++) {\n    $rowIndex = $dataRows[$r]\n    $rowValues = $newValues[$r]\n    for ($c = 1; $c -le $rowValues.Length; $c++) {\n        $t.Cell($rowIndex, $c).Range.Text = $rowValues[$c - 1]\n    }\n}\n"}
